# Refresh cryptos list (Price / Volume(1h) columns) with the latest scraped
# values. Price cells that look like a plain number (e.g. "309.45") are
# prefixed with a leading apostrophe so Excel stores them as text instead of
# auto-converting to a numeric value, matching the original inline-string
# cells (prices like "42.911.12" already contain multiple dots and can't be
# parsed as numbers, so they need no such prefix).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.911.12'
$ws.Range('E2').Value = '  +2.17%  '
$ws.Range('D3').Value = '2.296.17'
$ws.Range('E3').Value = '  +1.55%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''309.45'
$ws.Range('D6').Value = '''99.71'
$ws.Range('E6').Value = '  +4.51%  '
$ws.Range('E7').Value = '  +1.27%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = '''0.515'
$ws.Range('E9').Value = '  +5.48%  '
$ws.Range('D10').Value = '''35.96'
$ws.Range('E10').Value = '  +2.74%  '
$ws.Range('D11').Value = '''0.0817'
$ws.Range('E11').Value = '  +3.69%  '
$ws.Range('E12').Value = '  +0.74%  '
$ws.Range('D13').Value = '''7.13'
$ws.Range('E13').Value = '  +7.78%  '
$ws.Range('D14').Value = '2.656.96'
$ws.Range('E14').Value = '  +1.60%  '
$ws.Range('D15').Value = '''14.84'
$ws.Range('E15').Value = '  +3.43%  '
$ws.Range('D16').Value = '2.305.96'
$ws.Range('E16').Value = '  +3.47%  '
$ws.Range('D17').Value = '''0.799'
$ws.Range('E17').Value = '  +1.10%  '
$ws.Range('D18').Value = '42.877.41'
$ws.Range('E18').Value = '  +2.29%  '
$ws.Range('D19').Value = '''12.42'
$ws.Range('E19').Value = '  +0.48%  '
$ws.Range('D20').Value = '0.0₃0923'
$ws.Range('E20').Value = '  +2.40%  '
$ws.Range('D21').Value = '''6.05'
$ws.Range('E21').Value = '  +1.63%  '
$ws.Range('D22').Value = '''68.06'
$ws.Range('E22').Value = '  +0.66%  '
$ws.Range('D23').Value = '''239.29'
$ws.Range('E23').Value = '  +0.78%  '
$ws.Range('E24').Value = '  +4.24%  '
$ws.Range('E25').Value = '  +1.49%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').Value = '''24.05'
$ws.Range('E27').Value = '  +1.62%  '
$ws.Range('D28').Value = '''38.52'
$ws.Range('E28').Value = '  +5.33%  '
$ws.Range('D29').Value = '''9.62'
$ws.Range('E29').Value = '  +1.35%  '
$ws.Range('D30').Value = '''2.11'
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('D31').Value = '''168.01'
$ws.Range('E31').Value = '  +4.94%  '
$ws.Range('D32').Value = '''5.32'
$ws.Range('E32').Value = '  +2.02%  '
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('D34').Value = '''3.12'
$ws.Range('E34').Value = '  -1.78%  '
$ws.Range('E35').Value = '  +3.98%  '
$ws.Range('D36').Value = '''0.0736'
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('E37').Value = '  +0.22%  '
$ws.Range('E38').Value = '  +0.13%  '
$ws.Range('E39').Value = '  +1.60%  '
$ws.Range('E40').Value = '  +0.38%  '
$ws.Range('E41').Value = '  +6.01%  '
$ws.Range('D42').Value = '''2.28'
$ws.Range('E42').Value = '  -5.45%  '
$ws.Range('D43').Value = '1.963.62'
$ws.Range('E43').Value = '  -0.84%  '
$ws.Range('D44').Value = '''0.0287'
$ws.Range('E44').Value = '  +1.67%  '
$ws.Range('D45').Value = '''19.11'
$ws.Range('E45').Value = '  +1.03%  '
$ws.Range('D46').Value = '''2.99'
$ws.Range('E46').Value = '  +2.38%  '
$ws.Range('D47').Value = '''9.77'
$ws.Range('E47').Value = '  -1.34%  '
$ws.Range('D48').Value = '''2.96'
$ws.Range('E48').Value = '  +18.14%  '
$ws.Range('D49').Value = '''54.87'
$ws.Range('E49').Value = '  +3.28%  '
$ws.Range('D50').Value = '2.526.29'
$ws.Range('E50').Value = '  +1.58%  '
$ws.Range('D51').Value = '''1.54'
$ws.Range('E51').Value = '  +2.36%  '
